$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title "A" + " " + "slide" -> consolidated single run "A slide"
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = ""
$titleRange.Text = "A slide"

# Table cell "a" + " " + "table" -> consolidated single run "a table"
$tbl = $s.Shapes.Item(3).Table
$cellRange = $tbl.Cell(1, 2).Shape.TextFrame.TextRange
$cellRange.Text = ""
$cellRange.Text = "a table"

# TextBox "Plus" + " " + "an" + " " + "image" -> consolidated single run "Plus an image"
$boxRange = $s.Shapes.Item(7).TextFrame.TextRange
$boxRange.Text = ""
$boxRange.Text = "Plus an image"
